# Revised PCB to v1.2
# - Switched from linear regulator to DC-DC converter.
# - 1uF capacitors are now 10uF.
# - Added TAPR link to silkscreen and schematic (selection/active-sheet changes).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("pyboardbreakout")
$ws2 = $wb.Worksheets.Item("pyboardbreakout_data")

# --- Update source data on pyboardbreakout_data ---

# 40-pin WBUS connector -> 40-pin WBUS Connector (capitalization fix)
$ws2.Range("B7").Value = "40-pin WBUS Connector"

# 5V Linear Regulator -> 5V DC to DC Converter
$ws2.Range("B8").Value = "5V DC to DC Converter"
$ws2.Range("C8").Value = "R-78E5.0-1.0"
$ws2.Range("D8").Value = "dc_converter.pdf"

# 1uF 25V Capacitor (1206) -> 10uF 25V Capacitor (1206)
$ws2.Range("B11").Value = "10µF 25V Capacitor (1206)"

# Update the Order Link for the capacitor row: delete the old hyperlink that
# pointed at the 1uF part and add a new one pointing at the 10uF part.
$existingLinks = @($ws2.Hyperlinks)
foreach ($hl in $existingLinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$F$11') {
        $hl.Delete()
    }
}
$ws2.Range("F11").Value = "https://www.digikey.com/products/en?keywords=1276-1804-1-ND"
$ws2.Hyperlinks.Add($ws2.Range("F11"), "https://www.digikey.com/products/en?keywords=1276-1804-1-ND")

$wb.Application.Calculate()

# --- Update the active sheet / selections ---
$ws2.Range("C16").Select()
$ws1.Activate()
$ws1.Range("B8").Select()
